$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new "season record" columns ---
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy formatting (bold, centered, thin border) from an existing header
# cell so the new headers match the rest of row 1 exactly.
$ws.Range("AA1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# --- Data rows: season record repeated for every player row (2..45) ---
$wins = 70
$losses = 92
$ties = 0

for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 29).Value = $wins
    $ws.Cells.Item($r, 30).Value = $losses
    $ws.Cells.Item($r, 31).Value = $ties
}
